# Add a "Result" column to the credentials test-data sheet and record the
# outcome of the admin login test.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the existing header row's formatting into the new column (E) the
# same way the bordered/shaded header look is carried across D -> E, then
# write the new header text.
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Range("E1").Value = "Result"

# Record the test result for the first data row.
$ws.Range("E2").Value = "PASS"

$ws.Rows("1:1").RowHeight = 19.7

$ws.Range("E2").Select()
